$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Team member NIMs (column E, rows 1-2) ---
$ws.Range("E1").Value = "NIM : 201810370311073"
$ws.Range("E2").Value = "NIM : 201810370311076"

# --- Project title (C3) ---
$ws.Range("C3").Value = "Prosedur komputasi untuk pengenalan dan klasifikasi penyakit daun jagung dari daun sehat menggunakan metode CNN"

# --- Dataset title + link (C4, C5) ---
$ws.Range("C4").Value = "PlantVillage"
$ws.Range("C5").Value = "https://github.com/spMohanty/PlantVillage-Dataset"

# --- Article title (C6) - loses its border/style formatting in the source edit ---
$ws.Range("C6").ClearFormats()
$ws.Range("C6").Value = "Convolutional neural network for maize leaf disease image classification"

# --- Article link (C7) as a hyperlink ---
$ws.Range("C7").ClearFormats()
$ws.Hyperlinks.Add($ws.Range("C7"), "http://dx.doi.org/10.12928/telkomnika.v18i3.14840")

# --- Column E width ---
$ws.Columns("E").ColumnWidth = 70.14

# --- Page setup: A4 portrait ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection moves to E11 ---
$ws.Range("E11").Select() | Out-Null
